# corrected ICDC Breed 1-14 scripts
#
# The workbook has a single sheet ("startup") used as a small lookup table
# of tab names + the Neo4j/Excel queries used to populate each tab.
# This change corrects the "FilesTab" query (row 4, column B): it drops the
# now-redundant `File Type` and `Breed` projections from the RETURN clause.
# The row also gets a little shorter (less wrapped text -> smaller row
# height) and the workbook was left with the B4 cell selected/in view.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Corrected FilesTab query text (File Type + Breed lines removed) -----
$newFilesQuery = @'
MATCH (f:file)-->(parent)
WITH DISTINCT f, parent
MATCH (f)-[*]->(c:case)<--(demo:demographic)
WHERE demo.breed IN ['German Shepherd Dog'] 
OPTIONAL MATCH (s:study)<-[*]-(c)<--(diag:diagnosis)
OPTIONAL MATCH (samp:sample)-->(c)
WITH DISTINCT f, parent, c, demo, diag, s
RETURN  coalesce(f.file_name, '') AS `File Name`,
        coalesce(labels(parent)[0], '') AS `Association`,
        coalesce(f.file_description, '') AS `Description`,
        coalesce(f.file_format, '') AS `Format`,
        coalesce(f.file_size, '') AS `Size`,
        coalesce(c.case_id, '') AS `Case ID`,
         coalesce(diag.disease_term,'') AS Diagnosis , 
        coalesce(s.clinical_study_designation,'') AS `Study Code`
'@

$ws.Range("B4").Value = $newFilesQuery

# Row 4 now wraps to a shorter block of text.
$ws.Rows.Item(4).RowHeight = 217.5

# Leave the view with B4 selected/scrolled into place, as in the saved file.
$ws.Range("B4").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 4
$win.ScrollColumn = 1
